# Applies the "LetterFreq" workbook update:
#  - Adjust the window scroll position on the "High Damage Spells" and
#    "High Recoil Spells" sheets back to the default top-left cell.
#  - Add a new "Sheet1" worksheet at the end of the workbook that computes,
#    for every letter A-Z, the ratio of actual-to-expected frequency for
#    each of the six spell categories.

$wb = $excel.ActiveWorkbook

# --- Reset the stored scroll position on a couple of the existing sheets ---
$highDamage = $wb.Worksheets.Item("High Damage Spells")
$highDamage.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

$highRecoil = $wb.Worksheets.Item("High Recoil Spells")
$highRecoil.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

# --- Add the new "Sheet1" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet1"

# Header row
$newSheet.Range("A1").Value = "Letter:"
$newSheet.Range("B1").Value = "Damage"
$newSheet.Range("C1").Value = "Duration"
$newSheet.Range("D1").Value = "Instance"
$newSheet.Range("E1").Value = "AOE"
$newSheet.Range("F1").Value = "Knockback"
$newSheet.Range("G1").Value = "Recoil"

$categorySheets = @(
    "High Damage Spells",
    "High Duration Spells",
    "High Instance Spells",
    "AOE Spells",
    "High Knockback Spells",
    "High Recoil Spells"
)
$columns = @("B", "C", "D", "E", "F", "G")

$letters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = $i + 2
    $srcRow = $i + 1
    $newSheet.Range("A$row").Value = $letters[$i]

    for ($j = 0; $j -lt $categorySheets.Length; $j++) {
        $col = $columns[$j]
        $sheetName = $categorySheets[$j]
        $newSheet.Range("$col$row").Formula = "='$sheetName'!D$srcRow  / '$sheetName'!E$srcRow"
    }
}

# Make the new sheet the active one, with the same selection/scroll state
# recorded by the original author.
$newSheet.Activate()
$newSheet.Range("G28").Select()
